$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Insert a new column before column N (the 14th column). Excel shifts every
# column from N onward one position to the right (N->O, O->P, ... AR->AS),
# carrying each column's width/format along with it, and every existing
# cell reference shifts the same way (row spans, hyperlinks, etc.).
$ws.Columns("N").Insert()

# Populate the header for the freshly inserted column.
$ws.Range("N1").Value = "Another Amount"

# Populate the "amount" cell a few rows down; force it to stay literal text
# (rather than being auto-parsed as a currency number) the same way Excel
# treats an apostrophe-prefixed entry, after clearing the format the column
# insert implicitly inherited for this cell.
$ws.Range("N51").ClearFormats()
$ws.Range("N51").Value = "'$45.00"

# Reflect the saved view/selection state.
[void]$ws.Range("N34").Select()
